$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(7)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 1: "STUDENT NAME: APARNA S" -> "STUDENT NAME:  DHARHSINI T" ---
# Two-step set through a non-overlapping placeholder so the host doesn't keep a
# stale prefix run (it diffs old/new text and only fully replaces the run when
# there is no shared prefix/suffix).
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "`u{2060}1"
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "STUDENT NAME:  DHARHSINI T"

# --- Paragraph 2: "REGISTER NO :AND NMID: " -> "REGISTER NO :  2428CO469" ---
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "`u{2060}2"
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "REGISTER NO :  2428CO469"

# --- Paragraph 3: "DEPARTMENT: " -> "NMID: " + new run with the NMID value ---
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "`u{2060}3"
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "NMID: "
$para3 = $tr.Paragraphs(3, 1)
$nmidRun = $para3.InsertAfter("27280AA6F18645F999742B38ED971BD3")

# --- New paragraph inserted after paragraph 3: DEPARTMENT line ---
$para3 = $tr.Paragraphs(3, 1)
$deptPara = $para3.InsertAfter("`rDEPARTMENT: B.Sc.Artifical Intellegence & Machine Learning")

# --- College paragraph (now paragraph 5): "COLLEGE: COLLEGE/ UNIVERSITY" -> "COLLEGE: KPR College of Arts Science & Research" ---
$para5 = $tr.Paragraphs(5, 1)
$para5.Text = "`u{2060}5"
$para5 = $tr.Paragraphs(5, 1)
$para5.Text = "COLLEGE: KPR College of Arts Science & Research"

# --- Resize the textbox to match the extra line of text ---
$sh.Height = 181.65
